$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.070.67"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.103.53"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'579.15"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'172.66"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "'6.50"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "'36.66"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "3.616.38"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "67.007.06"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D18").Value = "3.101.90"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "'16.52"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'490.05"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "'0.700"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").Value = "'7.82"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").Value = "'83.69"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'13.05"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("D26").Value = "'10.52"
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "'28.19"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "0.0₃0929"
$ws.Range("E33").Value = "  -7.45%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").Value = "'46.93"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "'0.306"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").Value = "'384.31"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "2.796.73"
$ws.Range("E44").Value = "  -8.90%  "
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'135.01"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D48").Value = "'24.90"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -2.28%  "
